{"js": "// Applies the 2023 -> 2024 template refresh:\n//  1) blank-signature-date placeholder year bump (\"___ _____________ 2023 \u0440\u043e\u043a\u0443\" -> \"...2024...\")\n//  2) rector power-of-attorney clause: number 65 -> 70, add \"\u0432.\u043e.\" (acting) title, date 01.05.2023 -> 09.07.2024\n//  3) contract-number date stamp 2023 -> 2024\n//  4) total price 179 600,00 -> 199 600,00 (and the spelled-out words \"\u0441\u0442\u043e \u0441\u0456\u043c\u0434\u0435\u0441\u044f\u0442\" -> \"\u0441\u0442\u043e \u0434\u0435\u0432'\u044f\u043d\u043e\u0441\u0442\u043e\")\nconst body = context.document.body;\n\n// 1) Signature block blank date placeholder (table cell), e.g. \"  ___ _____________ 2023 \u0440\u043e\u043a\u0443\"\nconst placeholderDate = body.search(\"___ _____________ 2023 \u0440\u043e\u043a\u0443\", { matchCase: true });\nplaceholderDate.load(\"items\");\nawait context.sync();\nif (placeholderDate.items.length > 0) {\n  placeholderDate.items[0].insertText(\"___ _____________ 2024 \u0440\u043e\u043a\u0443\", Word.InsertLocation.replace);\n}\n\n// 2) Power-of-attorney clause for the rector who signs on behalf of the University\nconst attorneyClause = body.search(\"\u0414\u043e\u0440\u0443\u0447\u0435\u043d\u043d\u044f \u0440\u0435\u043a\u0442\u043e\u0440\u0430 \u2116 65 \u0432\u0456\u0434 01.05.2023\u0440\", { matchCase: true });\nattorneyClause.load(\"items\");\nawait context.sync();\nif (attorneyClause.items.length > 0) {\n  attorneyClause.items[0].insertText(\"\u0414\u043e\u0440\u0443\u0447\u0435\u043d\u043d\u044f \u0432.\u043e. \u0440\u0435\u043a\u0442\u043e\u0440\u0430 \u2116 70 \u0432\u0456\u0434 09.07.2024\u0440\", Word.InsertLocation.replace);\n}\n\n// 3) Contract number's date stamp (\"... \u2116 ______________ \u0432\u0456\u0434 ____________ 2023\u0440.\")\nconst contractDate = body.search(\"____________ 2023\u0440\", { matchCase: true });\ncontractDate.load(\"items\");\nawait context.sync();\nif (contractDate.items.length > 0) {\n  contractDate.items[0].insertText(\"____________ 2024\u0440\", Word.InsertLocation.replace);\n}\n\n// 4) Total tuition price, numeric and spelled-out Ukrainian words\nconst priceNumber = body.search(\"179\", { matchCase: true });\npriceNumber.load(\"text\");\nawait context.sync();\nfor (const item of priceNumber.items) {\n  if (item.text === \"179\") {\n    item.insertText(\"199\", Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n\nconst priceWords = body.search(\"\u0441\u0442\u043e \u0441\u0456\u043c\u0434\u0435\u0441\u044f\u0442 \u0434\u0435\u0432\\u2019\u044f\u0442\u044c \u0442\u0438\u0441\u044f\u0447 \u0448\u0456\u0441\u0442\u0441\u043e\u0442\", { matchCase: true });\npriceWords.load(\"items\");\nawait context.sync();\nif (priceWords.items.length > 0) {\n  priceWords.items[0].insertText(\"\u0441\u0442\u043e \u0434\u0435\u0432\\u2019\u044f\u043d\u043e\u0441\u0442\u043e \u0434\u0435\u0432\\u2019\u044f\u0442\u044c \u0442\u0438\u0441\u044f\u0447 \u0448\u0456\u0441\u0442\u0441\u043e\u0442\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Applies the 2023 -> 2024 template refresh:\n#  1) blank-signature-date placeholder year bump (\"___ _____________ 2023 \u0440\u043e\u043a\u0443\" -> \"...2024...\")\n#  2) rector power-of-attorney clause: number 65 -> 70, add \"\u0432.\u043e.\" (acting) title, date 01.05.2023 -> 09.07.2024\n#  3) contract-number date stamp 2023 -> 2024\n#  4) total price 179 600,00 -> 199 600,00 (and the spelled-out words \"\u0441\u0442\u043e \u0441\u0456\u043c\u0434\u0435\u0441\u044f\u0442\" -> \"\u0441\u0442\u043e \u0434\u0435\u0432'\u044f\u043d\u043e\u0441\u0442\u043e\")\n$d = $word.ActiveDocument\n\n# 1) Signature block blank date placeholder (table cell), e.g. \"  ___ _____________ 2023 \u0440\u043e\u043a\u0443\"\n$r1 = $d.Content\n$r1.Find.Execute(\"___ _____________ 2023 \u0440\u043e\u043a\u0443\", $false, $false, $false, $false, $false, $true, 1, $false, \"___ _____________ 2024 \u0440\u043e\u043a\u0443\", 2)\n\n# 2) Power-of-attorney clause for the rector who signs on behalf of the University\n$r2 = $d.Content\n$r2.Find.Execute(\"\u0414\u043e\u0440\u0443\u0447\u0435\u043d\u043d\u044f \u0440\u0435\u043a\u0442\u043e\u0440\u0430 \u2116 65 \u0432\u0456\u0434 01.05.2023\u0440\", $false, $false, $false, $false, $false, $true, 1, $false, \"\u0414\u043e\u0440\u0443\u0447\u0435\u043d\u043d\u044f \u0432.\u043e. \u0440\u0435\u043a\u0442\u043e\u0440\u0430 \u2116 70 \u0432\u0456\u0434 09.07.2024\u0440\", 2)\n\n# 3) Contract number's date stamp (\"... \u2116 ______________ \u0432\u0456\u0434 ____________ 2023\u0440.\")\n$r3 = $d.Content\n$r3.Find.Execute(\"____________ 2023\u0440\", $false, $false, $false, $false, $false, $true, 1, $false, \"____________ 2024\u0440\", 2)\n\n# 4) Total tuition price, numeric and spelled-out Ukrainian words\n$r4 = $d.Content\n$r4.Find.Execute(\"179\", $false, $true, $false, $false, $false, $true, 1, $false, \"199\", 2)\n\n$r5 = $d.Content\n$r5.Find.Execute(\"\u0441\u0442\u043e \u0441\u0456\u043c\u0434\u0435\u0441\u044f\u0442 \u0434\u0435\u0432\u2019\u044f\u0442\u044c \u0442\u0438\u0441\u044f\u0447 \u0448\u0456\u0441\u0442\u0441\u043e\u0442\", $false, $false, $false, $false, $false, $true, 1, $false, \"\u0441\u0442\u043e \u0434\u0435\u0432\u2019\u044f\u043d\u043e\u0441\u0442\u043e \u0434\u0435\u0432\u2019\u044f\u0442\u044c \u0442\u0438\u0441\u044f\u0447 \u0448\u0456\u0441\u0442\u0441\u043e\u0442\", 2)\n"}
